$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @(0.003323952356974284, 0.0004452725678422599, 0.001994689305623373, 0.0000008104673248279548)
    3  = @(0.003667672475179037, 0.0009570149351590361, 0.002659320831298828, 0.000939415976526479)
    4  = @(0.003346045811971029, 0.001229466340708883, 0.001662810643513997, 0.0004714277653650701)
    5  = @(0.002304315567016602, 0.0004865479313779051, 0.001662492752075195, 0.0004924637568278591)
    6  = @(0.002981265385945638, 0.00001439139168347066, 0.002005418141682943, 0.0000293918528865542)
    7  = @(0.002327283223470052, 0.0004694598580568356, 0.001340389251708984, 0.0004860375447421871)
    8  = @(0.004321177800496419, 0.001879861875499347, 0.002327521642049154, 0.0004708084360109519)
    9  = @(0.002992788950602213, 0.001410964096548703, 0.003010272979736328, 0.001437713415015756)
    10 = @(0.003342787424723307, 0.0004811230347141793, 0.002653121948242188, 0.0009279066708331331)
    11 = @(0.003325541814168294, 0.000474858483440153, 0.001329501469930013, 0.000444239793616943)
    12 = @(0.002660115559895833, 0.0004704916395046088, 0.00134126345316569, 0.0004618466899366043)
    13 = @(0.002671321233113607, 0.0004785059628636676, 0.001972357432047526, 0.0008025588762342157)
    14 = @(0.002657492955525716, 0.0009384700085535198, 0.002316554387410482, 0.0004552983971299241)
    15 = @(0.004008769989013672, 0.001434849717396358, 0.001972277959187826, 0.001416918524317377)
    16 = @(0.003325939178466797, 0.001245026035683998, 0.002005736033121745, 0.0008270858453644146)
    17 = @(0.002326250076293945, 0.0004703591919349732, 0.001651763916015625, 0.0004627821924254846)
}

foreach ($row in $values.Keys) {
    $vals = $values[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
}
